$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats to reflect the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.03      # Total P&L $
$summary.Range("B5").Value = 0.01      # Total P&L %
$summary.Range("B6").Value = 102       # Total Trades
$summary.Range("B8").Value = 42        # Losing Trades
$summary.Range("B9").Value = 41.18     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 102        # Trades
$status.Range("E4").Value = 0.03       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 41.18      # Win Rate %

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append new trade #102 as row 103
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A103").Value = 102

    # "2026-02-17" parses as a date literal via Value, so force text storage
    # then restore the default (Normal) style so no stray format is left on
    # the cell.
    $ws.Range("B103").NumberFormat = "@"
    $ws.Range("B103").Value = "2026-02-17"
    $ws.Range("B103").Style = "Normal"

    $ws.Range("C103").Value = "09:18:03"
    $ws.Range("D103").Value = "MarketMaking"
    $ws.Range("E103").Value = "DOWN"
    $ws.Range("F103").Value = 0.879837
    $ws.Range("G103").Value = 0.8
    $ws.Range("H103").Value = "CLOSED"
    $ws.Range("I103").Value = -9.074
    $ws.Range("J103").Value = -0.08
    $ws.Range("K103").Value = 100.02
    $ws.Range("L103").Value = 0
    $ws.Range("M103").Value = 0
    $ws.Range("N103").Value = 0.6
    $ws.Range("O103").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P103").Value = "early_exit"
    $ws.Range("Q103").Value = 0.13
}
